{"js": "// Commit: \"add materials for week2\"\n// The \"Examples\" table's Output column contains a cell whose single run\n// reads \"4.00\" (the printed interest amount). The authoritative edit\n// re-keys that value to \"0.40\", but does so by exploding the text into\n// four individual single-character runs (\"0\", \".\", \"4\", \"0\") that each\n// carry the same run formatting (Consolas, bCs, noProof) as the\n// original run. Reproduce that exact run layout.\n\nconst results = context.document.body.search(\"4.00\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find target text '4.00' in the document body.\");\n}\n\nconst target = results.items[0];\n\n// Run properties shared by every run in the cell (Consolas font, bold\n// complex-script flag off, noProof) \u2014 identical to the original run's\n// <w:rPr>.\nconst rPr =\n  '<w:rPr><w:rFonts w:ascii=\"Consolas\" w:hAnsi=\"Consolas\"/><w:bCs/><w:noProof/></w:rPr>';\n\n// Rebuild the run as four sibling runs inside the same paragraph, each\n// holding one character of \"0.40\". Using insertOoxml(\"Replace\") lets us\n// control the run boundaries precisely (plain insertText would collapse\n// the text back into a single run). The first run keeps the original\n// w:rsidRPr attribute; the three newly-introduced runs do not carry it,\n// matching the diff.\nconst ooxml =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p>' +\n  '<w:r w:rsidRPr=\"002045C4\">' + rPr + '<w:t>0</w:t></w:r>' +\n  '<w:r>' + rPr + '<w:t>.</w:t></w:r>' +\n  '<w:r>' + rPr + '<w:t>4</w:t></w:r>' +\n  '<w:r>' + rPr + '<w:t>0</w:t></w:r>' +\n  '</w:p>' +\n  '</w:body></w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>';\n\ntarget.insertOoxml(ooxml, \"Replace\");\nawait context.sync();\n", "ps1": "# Commit: \"add materials for week2\"\n#\n# The \"Examples\" table's Output column holds a cell whose single run\n# reads \"4.00\" (a printed interest amount). The authoritative edit\n# re-keys that value to \"0.40\", exploding the text into four\n# individual single-character runs (\"0\", \".\", \"4\", \"0\") that each carry\n# the same run formatting (Consolas, bCs, noProof) the original run\n# had. We reproduce that exact run layout below.\n\n$d = $word.ActiveDocument\n$rng = $d.Content\n\n$found = $rng.Find.Execute(\"4.00\")\nif (-not $found) {\n    throw \"Could not find target text '4.00' in the document.\"\n}\n\n# $rng now covers exactly the four characters \"4.00\" inside its\n# paragraph. Word's InsertXML replaces the whole paragraph that the\n# range lives in whenever the payload carries a <w:p> element, so the\n# fragment below restates that paragraph's own (unchanged) identity/\n# formatting attributes exactly as they already are in the document,\n# and then supplies four sibling runs in place of the single original\n# run. The first run keeps the original run's w:rsidRPr attribute; the\n# three newly-introduced runs do not carry it -- matching the diff.\n$xml = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' + `\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>' + `\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' + `\n  '<w:body>' + `\n  '<w:p w14:paraId=\"6D8C5E05\" w14:textId=\"77777777\" w:rsidR=\"002A3A59\" w:rsidRPr=\"002045C4\" w:rsidRDefault=\"002A3A59\" w:rsidP=\"00F06727\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\">' + `\n  '<w:pPr><w:spacing w:before=\"0\" w:after=\"0\"/><w:rPr><w:rFonts w:ascii=\"Consolas\" w:hAnsi=\"Consolas\"/><w:bCs/><w:noProof/></w:rPr></w:pPr>' + `\n  '<w:r w:rsidRPr=\"002045C4\"><w:rPr><w:rFonts w:ascii=\"Consolas\" w:hAnsi=\"Consolas\"/><w:bCs/><w:noProof/></w:rPr><w:t>0</w:t></w:r>' + `\n  '<w:r><w:rPr><w:rFonts w:ascii=\"Consolas\" w:hAnsi=\"Consolas\"/><w:bCs/><w:noProof/></w:rPr><w:t>.</w:t></w:r>' + `\n  '<w:r><w:rPr><w:rFonts w:ascii=\"Consolas\" w:hAnsi=\"Consolas\"/><w:bCs/><w:noProof/></w:rPr><w:t>4</w:t></w:r>' + `\n  '<w:r><w:rPr><w:rFonts w:ascii=\"Consolas\" w:hAnsi=\"Consolas\"/><w:bCs/><w:noProof/></w:rPr><w:t>0</w:t></w:r>' + `\n  '</w:p>' + `\n  '</w:body></w:document>' + `\n  '</pkg:xmlData></pkg:part></pkg:package>'\n\n$rng.InsertXML($xml) | Out-Null\n"}
